$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45309

$ws.Range("D14").Value = 38.22
$ws.Range("D15").Value = 43.55
$ws.Range("D16").Value = 52.65
$ws.Range("D17").Value = 66.95
$ws.Range("D18").Value = 133.9
$ws.Range("D19").Value = 201.5
$ws.Range("D20").Value = 314.6
$ws.Range("D21").Value = 432.9
